$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the redundant data point: row 9 (Drawdown Start 43126 / 9 days / text "-10.20%" / Conclusion 43139)
# Deleting the entire row shifts rows 10 and 11 up to become rows 9 and 10.
$ws.Rows(9).Delete()

# That row's "C" cell was the only one using the "Percent" named cell style
# (and the placeholder shared string "-10.20%"); with it gone, drop the now
# unused style definition too.
$wb.Styles("Percent").Delete()

# Leftover selection position recorded in the saved file.
$ws.Range("C16").Select()
